$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Checklist")

# Delete old row 18 (TestData_07 entry); rows 19-24 shift up to 18-23
$ws.Rows.Item(18).Delete()

# Row 18 updates
$ws.Range("C18").Value = 87
$ws.Range("E18").Value = "X"
$ws.Range("O18").Value = "HOLA"
$ws.Range("P18").Value = 7
$ws.Range("Q18").Value = "JULIO"
$ws.Range("S18").Value = 5

# Row 19 updates
$ws.Range("C19").Value = 4
$ws.Range("E19").Value = "X"
$ws.Range("F19").Value = "X"
$ws.Range("N19").Value = "BB96"
$ws.Range("O19").Value = "ADIOS"
$ws.Range("P19").Value = 6
$ws.Range("R19").Value = 7
$ws.Range("S19").Value = 3

# Row 20 updates
$ws.Range("C20").Value = 9
$ws.Range("D20").Value = "X"
$ws.Range("E20").Value = "X"
$ws.Range("N20").Value = "BB75"
$ws.Range("O20").Value = ":("
$ws.Range("Q20").Value = "JULIO Y RUBEN"
$ws.Range("R20").Value = 6
$ws.Range("S20").Value = 6

# Row 21 updates
$ws.Range("D21").Value = "X"
$ws.Range("F21").Value = "X"
$ws.Range("M21").Value = "datablock"
$ws.Range("N21").Value = "BB89"
$ws.Range("O21").ClearContents()
$ws.Range("P21").Value = 10
$ws.Range("Q21").Value = "JEJ"
$ws.Range("R21").Value = 5
$ws.Range("S21").Value = 9

# Row 24 updates
$ws.Range("A24").Value = "TestData_07"
$ws.Range("B24").Value = 123456
$ws.Range("C24").Value = 91
$ws.Range("F24").Value = "X"
$ws.Range("G24").Value = "X"
$ws.Range("H24").Value = "X"
$ws.Range("L24").Value = "X"
$ws.Range("M24").Value = "ee_range"
$ws.Range("O24").Value = "description=- Component: Test`n- REPROG info: To be evaluated."
$ws.Range("S24").Value = 9
